# Update visitor-count ("想去人数") figures in the workbook.
# Sheet "展览" (sheet1) and "全部类型" (sheet4) share the same first 21
# data rows; "全部类型" additionally carries the "演出" row (as its row 22)
# and the last "展览" row (as its row 23). "演出" (sheet2) holds that same
# row independently as its own row 2.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# Column F updates shared by 展览 (rows 2-22) and 全部类型 (rows 2-21, same values)
$updates = @{
    2  = 1077
    3  = 368
    4  = 1478
    5  = 8702
    6  = 86
    7  = 490
    8  = 642
    9  = 278
    11 = 10
    12 = 3546
    14 = 363
    15 = 76
    16 = 1151
    18 = 1112
    19 = 306
    21 = 2287
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $wsExpo.Cells.Item($row, 6).Value = $value
    $wsAll.Cells.Item($row, 6).Value = $value
}

# 展览 row 22 is its own last row (F22: 50 -> 52); 全部类型 carries this
# same record shifted down to row 23 (F23: 50 -> 52).
$wsExpo.Cells.Item(22, 6).Value = 52
$wsAll.Cells.Item(23, 6).Value = 52

# 演出 row 2 (F2: 34 -> 35); 全部类型 carries this same record as row 22.
$wsShow.Cells.Item(2, 6).Value = 35
$wsAll.Cells.Item(22, 6).Value = 35
